$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(4)  # template: strategy_id-5009

$newNames = @("strategy_id-5011", "strategy_id-5012", "strategy_id-5013", "strategy_id-5014", "strategy_id-5015")

$allValues = @()
$allValues += ,@(0.9829063213546697, 0.9658126427093395, 0.9487189640640092, 0.931625285418679, 0.9145316067733488, 0.8974379281280185, 0.8803442494826883, 0.863250570837358, 0.8461568921920278, 0.8290632135466975, 0.8119695349013674, 0.7948758562560371, 0.7777821776107068, 0.7606884989653766, 0.7435948203200462, 0.726501141674716, 0.7094074630293857, 0.6923137843840557, 0.6752201057387253, 0.6581264270933951, 0.6410327484480649, 0.6239390698027346, 0.6068453911574043, 0.5897517125120741, 0.5726580338667439)
$allValues += ,@(0.9862304763364699, 0.9724609526729397, 0.9586914290094095, 0.9449219053458794, 0.9311523816823493, 0.9173828580188191, 0.903613334355289, 0.8898438106917588, 0.8760742870282288, 0.8623047633646985, 0.8485352397011684, 0.8347657160376383, 0.820996192374108, 0.8072266687105779, 0.7934571450470478, 0.7796876213835175, 0.7659180977199874, 0.7521485740564573, 0.7383790503929272, 0.724609526729397, 0.7108400030658669, 0.6970704794023368, 0.6833009557388066, 0.6695314320752765, 0.6557619084117463)
$allValues += ,@(0.9831702487644064, 0.9663404975288129, 0.9495107462932192, 0.9326809950576255, 0.915851243822032, 0.8990214925864384, 0.8821917413508448, 0.8653619901152512, 0.8485322388796576, 0.831702487644064, 0.8148727364084704, 0.7980429851728768, 0.7812132339372831, 0.7643834827016895, 0.7475537314660958, 0.7307239802305022, 0.7138942289949086, 0.6970644777593151, 0.6802347265237214, 0.6634049752881279, 0.6465752240525343, 0.6297454728169407, 0.6129157215813471, 0.5960859703457535, 0.5792562191101598)
$allValues += ,@(0.9888585559026216, 0.9777171118052433, 0.966575667707865, 0.9554342236104866, 0.9442927795131083, 0.9331513354157299, 0.9220098913183516, 0.9108684472209733, 0.899727003123595, 0.8885855590262166, 0.8774441149288383, 0.86630267083146, 0.8551612267340816, 0.8440197826367033, 0.8328783385393249, 0.8217368944419465, 0.8105954503445683, 0.79945400624719, 0.7883125621498116, 0.7771711180524332, 0.766029673955055, 0.7548882298576767, 0.7437467857602983, 0.7326053416629199, 0.7214638975655416)
$allValues += ,@(0.9813124864730004, 0.9626249729460008, 0.9439374594190011, 0.9252499458920015, 0.906562432365002, 0.8878749188380024, 0.8691874053110027, 0.8504998917840031, 0.8318123782570036, 0.813124864730004, 0.7944373512030043, 0.7757498376760047, 0.7570623241490051, 0.7383748106220054, 0.7196872970950059, 0.7009997835680062, 0.6823122700410066, 0.663624756514007, 0.6449372429870075, 0.6262497294600078, 0.6075622159330083, 0.5888747024060086, 0.5701871888790091, 0.5514996753520094, 0.5328121618250098)

for ($i = 0; $i -lt $newNames.Count; $i++) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $src.Cells.Copy($new.Range("A1"))
    $new.Name = $newNames[$i]
    $rowVals = $allValues[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $col = 21 + $j  # column U is index 21
        $new.Cells.Item(2, $col).Value2 = $rowVals[$j]
    }
}

$wb.Worksheets.Item(1).Select()
